# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on the zh-cn and de-de
# report sheets to reflect the latest handback run timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 01:03:28"
$wsZhCn.Range("H2").Value = "2016-03-22 01:03:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 01:03:32"
$wsDeDe.Range("H2").Value = "2016-03-22 01:03:56"
